# Update "想去人数" (F column) counts across sheets, per upstream data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - F column updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 514
$ws1.Range("F6").Value = 459
$ws1.Range("F7").Value = 484
$ws1.Range("F8").Value = 290
$ws1.Range("F9").Value = 15
$ws1.Range("F10").Value = 13005
$ws1.Range("F11").Value = 13005
$ws1.Range("F16").Value = 208
$ws1.Range("F17").Value = 173
$ws1.Range("F18").Value = 218
$ws1.Range("F19").Value = 2798
$ws1.Range("F22").Value = 2137
$ws1.Range("F23").Value = 188
$ws1.Range("F27").Value = 2490
$ws1.Range("F29").Value = 1207
$ws1.Range("F30").Value = 4447
$ws1.Range("F32").Value = 4076
$ws1.Range("F33").Value = 1116
$ws1.Range("F34").Value = 2730
$ws1.Range("F35").Value = 3134
$ws1.Range("F36").Value = 116
$ws1.Range("F37").Value = 1453
$ws1.Range("F39").Value = 808
$ws1.Range("F40").Value = 74
$ws1.Range("F41").Value = 206
$ws1.Range("F42").Value = 737
$ws1.Range("F43").Value = 1161
$ws1.Range("F44").Value = 103
$ws1.Range("F45").Value = 199
$ws1.Range("F46").Value = 475
$ws1.Range("F47").Value = 138
$ws1.Range("F48").Value = 248
$ws1.Range("F49").Value = 290

# Sheet "演出" (sheet2) - F column updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 80
$ws2.Range("F11").Value = 191
$ws2.Range("F13").Value = 20

# Sheet "全部类型" (sheet4) - F column updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 514
$ws4.Range("F5").Value = 459
$ws4.Range("F6").Value = 484
$ws4.Range("F7").Value = 290
$ws4.Range("F8").Value = 13005
$ws4.Range("F11").Value = 80
$ws4.Range("F13").Value = 208
$ws4.Range("F14").Value = 173
$ws4.Range("F16").Value = 218
$ws4.Range("F17").Value = 2798
$ws4.Range("F18").Value = 2137
$ws4.Range("F19").Value = 188
$ws4.Range("F24").Value = 2490
$ws4.Range("F25").Value = 1207
$ws4.Range("F26").Value = 191
$ws4.Range("F27").Value = 20
$ws4.Range("F28").Value = 4447
$ws4.Range("F30").Value = 4077
$ws4.Range("F31").Value = 1117
$ws4.Range("F32").Value = 2730
$ws4.Range("F33").Value = 3134
$ws4.Range("F34").Value = 116
$ws4.Range("F36").Value = 1453
$ws4.Range("F39").Value = 808
$ws4.Range("F40").Value = 74
$ws4.Range("F41").Value = 206
$ws4.Range("F42").Value = 737
$ws4.Range("F44").Value = 1161
$ws4.Range("F45").Value = 103
$ws4.Range("F46").Value = 199
$ws4.Range("F47").Value = 475
$ws4.Range("F48").Value = 138
$ws4.Range("F49").Value = 248
$ws4.Range("F50").Value = 290
